# GL90 - Debit Credit Save Close changes
# Applies the data edits made to the "Input_Value" sheet:
#  - Clears the EnterJournalName value in A2
#  - Updates DebitAmt (F2) and CreditAmt (G2) from 275 to 511
#  - Moves the active cell selection to E11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Clear the journal name entry in A2
$ws.Range("A2").Value = ""

# Update Debit/Credit amounts (stored as text, column is formatted as Text "@")
$ws.Range("F2").Value = "511"
$ws.Range("G2").Value = "511"

# Update the saved selection / active cell
[void]$ws.Range("E11").Select()
